# Generate Report for Handoff
# - Replace the old source file (dcfdb086-89ff-448d-b135-7d8d584c8f47.md) references
#   with the new one (4470401b-bb0b-45c7-9139-9310456c4006.md)
# - Remove the "Handoff transform failed" row (a2035c89-4d62-4576-a058-71adbf6cabaa.md)
#   from every sheet
# - Refresh the handoff package hash and timestamps for both locales

$wb = $excel.ActiveWorkbook

$oldMd = "dcfdb086-89ff-448d-b135-7d8d584c8f47.md"
$newMd = "4470401b-bb0b-45c7-9139-9310456c4006.md"

$newZhXlf = "4470401b-bb0b-45c7-9139-9310456c4006.b50252d8a42fc8c2bc2cc0fa86f1dcc73b54a5ba.zh-cn.xlf"
$newDeXlf = "4470401b-bb0b-45c7-9139-9310456c4006.b50252d8a42fc8c2bc2cc0fa86f1dcc73b54a5ba.de-de.xlf"

$newZhTime = "2016-02-18 03:58:21"
$newDeTime = "2016-02-18 03:58:34"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/1dac61efeea7626c3a2b5a254a988728327d3a3c/e2e/$newMd"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/1dac61efeea7626c3a2b5a254a988728327d3a3c/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fddeb1410fba48f3988bd877666f8d6badb1d26/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newZhXlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/55d8ad65ff8dc7d5e47d88aeaabd1742d924854a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newDeXlf"

$hlColor = 15570276  # matches the workbook's existing hyperlink font color (FF6495ED)

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): File Name / zh-cn / de-de summary
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Update the values first (while the old hyperlinks still exist) so the
# existing cell formatting/style is preserved.
$ws1.Range("A2").Value = $newMd

# Drop all hyperlinks, then remove the obsolete "Handoff transform failed" row
# (this shifts the ".localization-config" row up from row 4 to row 3).
$ws1.Hyperlinks.Delete()
$ws1.Rows.Item(3).Delete()

# Re-create the hyperlinks that remain.
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

$ws1.Range("A2:A3").Font.Underline = $true
$ws1.Range("A2:A3").Font.Color = $hlColor

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): handoff details for the zh-cn locale
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newMd
$ws2.Range("C2").Value = $newZhXlf
$ws2.Range("D2").Value = $newZhTime

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhXlfUrl, "", "", $newZhXlf) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

$ws2.Range("A2").Font.Underline = $true
$ws2.Range("A2").Font.Color = $hlColor
$ws2.Range("C2").Font.Underline = $true
$ws2.Range("C2").Font.Color = $hlColor
$ws2.Range("A3").Font.Underline = $true
$ws2.Range("A3").Font.Color = $hlColor

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3): handoff details for the de-de locale
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newMd
$ws3.Range("C2").Value = $newDeXlf
$ws3.Range("D2").Value = $newDeTime

$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deXlfUrl, "", "", $newDeXlf) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

$ws3.Range("A2").Font.Underline = $true
$ws3.Range("A2").Font.Color = $hlColor
$ws3.Range("C2").Font.Underline = $true
$ws3.Range("C2").Font.Color = $hlColor
$ws3.Range("A3").Font.Underline = $true
$ws3.Range("A3").Font.Color = $hlColor
